# Split "natural gas nonpeaker" into "natural gas steam turbine" and
# "natural gas combined cycle" on the ETS sheet.
#
# Row 3 currently holds "natural gas nonpeaker" (value 1 for every year,
# columns B:AF). We rename row 3 to "natural gas steam turbine" (keeping its
# values), insert a new row 4 for "natural gas combined cycle" (copying the
# same values), and let every row below shift down by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Insert a new blank row before the current row 4 ("nuclear"), pushing
# "nuclear" and everything beneath it down by one row.
$ws.Rows(4).Insert()

# Row 3: rename "natural gas nonpeaker" -> "natural gas steam turbine".
# Its numeric shareweights (1 for every year) are unchanged.
$ws.Range("A3").Value = "natural gas steam turbine"

# Row 4 (new): "natural gas combined cycle", with the same shareweights
# that "natural gas nonpeaker" used to have (1 for every year, 2020-2050).
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("B4:AF4").Value = 1
